$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "91.983.81"
Set-TextValue "E2" "  -3.13%  "
Set-TextValue "D3" "3.278.65"
Set-TextValue "E3" "  -5.29%  "
Set-TextValue "E4" "  +0.19%  "
Set-TextValue "D5" "226.35"
Set-TextValue "E5" "  -5.61%  "
Set-TextValue "D6" "605.81"
Set-TextValue "E6" "  -6.00%  "
Set-TextValue "D7" "1.34"
Set-TextValue "E7" "  -8.63%  "
Set-TextValue "E8" "  -7.52%  "
Set-TextValue "E9" "  +0.07%  "
Set-TextValue "D10" "0.919"
Set-TextValue "E10" "  -9.21%  "
Set-TextValue "D11" "3.275.65"
Set-TextValue "E11" "  -5.39%  "
Set-TextValue "D12" "41.10"
Set-TextValue "E12" "  -1.88%  "
Set-TextValue "D13" "0.190"
Set-TextValue "E13" "  -4.05%  "
Set-TextValue "E14" "  -4.35%  "
Set-TextValue "D15" "91.844.21"
Set-TextValue "E15" "  -3.01%  "
Set-TextValue "D16" "3.886.01"
Set-TextValue "E16" "  -5.35%  "
Set-TextValue "E17" "  -7.03%  "
Set-TextValue "D18" "7.91"
Set-TextValue "E18" "  -7.19%  "
Set-TextValue "D19" "3.276.29"
Set-TextValue "E19" "  -5.12%  "
Set-TextValue "D20" "16.88"
Set-TextValue "E20" "  -5.85%  "
Set-TextValue "D21" "10.54"
Set-TextValue "E21" "  -7.94%  "
Set-TextValue "E22" "  +5.70%  "
Set-TextValue "D23" "479.88"
Set-TextValue "E23" "  -4.71%  "
Set-TextValue "D24" "0.433"
Set-TextValue "E24" "  -15.26%  "
Set-TextValue "D25" "0.0000175"
Set-TextValue "E25" "  -9.52%  "
Set-TextValue "D26" "5.96"
Set-TextValue "E26" "  -8.35%  "
Set-TextValue "D27" "88.43"
Set-TextValue "E27" "  -3.77%  "
Set-TextValue "D28" "11.53"
Set-TextValue "E28" "  -5.15%  "
Set-TextValue "D29" "3.461.64"
Set-TextValue "E29" "  -4.84%  "
Set-TextValue "E30" "  +0.06%  "
Set-TextValue "E31" "  -7.97%  "
Set-TextValue "E32" "  -3.10%  "
Set-TextValue "B33" "Binance-PegBSC-USD"
Set-TextValue "C33" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  -0.27%  "
Set-TextValue "B34" "PancakeSwap"
Set-TextValue "C34" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D34" "2.55"
Set-TextValue "E34" "  -7.13%  "
Set-TextValue "D35" "0.169"
Set-TextValue "E35" "  -8.04%  "
Set-TextValue "D36" "27.61"
Set-TextValue "E36" "  -10.98%  "
Set-TextValue "D37" "0.515"
Set-TextValue "E37" "  -9.18%  "
Set-TextValue "D38" "534.18"
Set-TextValue "E38" "  +1.68%  "
Set-TextValue "E39" "  -0.05%  "
Set-TextValue "D40" "7.17"
Set-TextValue "E40" "  -7.40%  "
Set-TextValue "E41" "  -3.81%  "
Set-TextValue "E42" "  -7.92%  "
Set-TextValue "E43" "  -9.03%  "
Set-TextValue "D44" "23.79"
Set-TextValue "E44" "  -1.21%  "
Set-TextValue "E45" "  -3.83%  "
Set-TextValue "D46" "3.54"
Set-TextValue "E46" "  +1.28%  "
Set-TextValue "E47" "  -4.20%  "
Set-TextValue "E48" "  -8.54%  "
Set-TextValue "D49" "51.25"
Set-TextValue "E49" "  -3.99%  "
Set-TextValue "E50" "  -5.50%  "
Set-TextValue "E51" "  -4.00%  "
